$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add formulas computing ratio B/A for rows 4 and 5
$ws.Range("C4").Formula = "=B4/A4"
$ws.Range("C5").Formula = "=B5/A5"

# Update the active cell selection
$ws.Range("K14").Select()
